$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values to re-pull/push the repulled data & mean calc
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -1
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 2
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 7
$ws.Range("F22").Value = 0
$ws.Range("F33").Value = -2
